# treasureManager & prepare rebuild storeItem structure
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update storeItem "can sell" / related flags in column D (and one E value)
$ws.Range("D16").Value = 1
$ws.Range("D17").Value = 1

$ws.Range("D28").Value = 1
$ws.Range("D29").Value = 0
$ws.Range("D30").Value = 1
$ws.Range("D32").Value = 1

$ws.Range("D35").Value = 1
$ws.Range("D36").Value = 1
$ws.Range("D37").Value = 1
$ws.Range("D38").Value = 1
$ws.Range("D39").Value = 1

$ws.Range("D40").Value = 0
$ws.Range("D41").Value = 0
$ws.Range("E41").Value = 16
$ws.Range("D42").Value = 0
$ws.Range("D43").Value = 0
$ws.Range("D44").Value = 1

# Move the active selection (was G34, scrolled to A20) to F12
$ws.Range("F12").Select()
